$d = $word.ActiveDocument
$d.Content.Find.Execute("have to thank Phil again.", $true, $false, $false, $false, $false, $true, 1, $false, "TESTREPLACE", 2)
